# Apodeixi "posting label" sheet: insert a new "Posting Label" worksheet in
# front of the existing "Manifest" worksheet, leaving the Manifest data
# completely untouched.

$wb = $excel.ActiveWorkbook
$manifest = $wb.Worksheets.Item("Manifest")

# --- Create the new sheet and place it immediately before "Manifest" ---
$postingLabel = $wb.Worksheets.Add()
$postingLabel.Name = "Posting Label"
$postingLabel.Move($manifest)

# --- Column widths (B ~5.71 chars, C ~19.71 chars, matching Manifest's units) ---
$postingLabel.Columns.Item(2).ColumnWidth = 4.833333333333333
$postingLabel.Columns.Item(3).ColumnWidth = 18.833333333333332

# --- B2: "TBD" -- header-style cell: bold white text on a gray fill,
#     centered/top aligned, wrapped, thin white right border (same look as
#     the Manifest header row). ---
$b2 = $postingLabel.Range("B2")
$b2.Value = "TBD"
$b2.Interior.Color = 8421504
$b2.Font.Bold = $true
$b2.Font.Color = 16777215
$b2.HorizontalAlignment = -4108
$b2.VerticalAlignment = -4160
$b2.WrapText = $true
$b2.Borders.Item(10).LineStyle = 1
$b2.Borders.Item(10).Color = 16777215

# --- C2: "In test_as_excel.py" -- body-style cell: light gray fill,
#     top aligned, wrapped (same look as Manifest's derived-value cells). ---
$c2 = $postingLabel.Range("C2")
$c2.Value = "In test_as_excel.py"
$c2.Interior.Color = 15921906
$c2.VerticalAlignment = -4160
$c2.WrapText = $true

# --- Protect the new sheet (matches Manifest's protection) ---
$postingLabel.Protect()

# --- Make "Posting Label" the selected/active tab ---
$postingLabel.Activate()

Write-Host "Inserted 'Posting Label' sheet before 'Manifest'."
